# Updates cryptos list values (Coin/Link/Price/Volume) per the Nov 8 2024 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '75.936.31'
$cell.ClearFormats()
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.51%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.908.08'
$cell.ClearFormats()
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.36%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.08%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '198.61'
$cell.ClearFormats()
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +5.59%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '599.61'
$cell.ClearFormats()
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.45%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.09%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.549'
$cell.ClearFormats()
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.34%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.200'
$cell.ClearFormats()
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +4.09%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.906.04'
$cell.ClearFormats()
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.35%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.432'
$cell.ClearFormats()
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +17.00%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.161'
$cell.ClearFormats()
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.20%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.88'
$cell.ClearFormats()
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.90%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.433.67'
$cell.ClearFormats()
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.20%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '75.694.86'
$cell.ClearFormats()
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.17%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0000192'
$cell.ClearFormats()
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.42%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '27.33'
$cell.ClearFormats()
$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.21%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.906.85'
$cell.ClearFormats()
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +3.34%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.98'
$cell.ClearFormats()
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +5.42%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.75'
$cell.ClearFormats()
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.19%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '370.59'
$cell.ClearFormats()
$cell = $ws.Cells.Item(21, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.06%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.30'
$cell.ClearFormats()
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +2.28%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.29'
$cell.ClearFormats()
$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +4.79%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.06%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '71.02'
$cell.ClearFormats()
$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.60%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(26, 2)
$cell.NumberFormat = "@"
$cell.Value = 'NEARProtocol'
$cell.ClearFormats()
$cell = $ws.Cells.Item(26, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell.ClearFormats()
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.20'
$cell.ClearFormats()
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.19%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(27, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Aptos'
$cell.ClearFormats()
$cell = $ws.Cells.Item(27, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell.ClearFormats()
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.65'
$cell.ClearFormats()
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.82%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(28, 2)
$cell.NumberFormat = "@"
$cell.Value = 'PEPE'
$cell.ClearFormats()
$cell = $ws.Cells.Item(28, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$cell.ClearFormats()
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0000108'
$cell.ClearFormats()
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +5.59%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(29, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Binance-PegBSC-USD'
$cell.ClearFormats()
$cell = $ws.Cells.Item(29, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$cell.ClearFormats()
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.997'
$cell.ClearFormats()
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.40%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(30, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Fetch.AI'
$cell.ClearFormats()
$cell = $ws.Cells.Item(30, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell.ClearFormats()
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.41'
$cell.ClearFormats()
$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +2.17%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(31, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Bittensor'
$cell.ClearFormats()
$cell = $ws.Cells.Item(31, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell.ClearFormats()
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '503.80'
$cell.ClearFormats()
$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.05%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(32, 2)
$cell.NumberFormat = "@"
$cell.Value = 'InternetComputer(DFINITY)'
$cell.ClearFormats()
$cell = $ws.Cells.Item(32, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell.ClearFormats()
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.70'
$cell.ClearFormats()
$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.97%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(33, 2)
$cell.NumberFormat = "@"
$cell.Value = 'PancakeSwap'
$cell.ClearFormats()
$cell = $ws.Cells.Item(33, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell.ClearFormats()
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.81'
$cell.ClearFormats()
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.14%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(34, 2)
$cell.NumberFormat = "@"
$cell.Value = 'FirstDigitalUSD'
$cell.ClearFormats()
$cell = $ws.Cells.Item(34, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$cell.ClearFormats()
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.ClearFormats()
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.04%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(35, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Monero'
$cell.ClearFormats()
$cell = $ws.Cells.Item(35, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell.ClearFormats()
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '165.06'
$cell.ClearFormats()
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.24%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(36, 2)
$cell.NumberFormat = "@"
$cell.Value = 'EthereumClassic'
$cell.ClearFormats()
$cell = $ws.Cells.Item(36, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell.ClearFormats()
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '20.15'
$cell.ClearFormats()
$cell = $ws.Cells.Item(36, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.38%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(37, 2)
$cell.NumberFormat = "@"
$cell.Value = 'WhiteBITCoin'
$cell.ClearFormats()
$cell = $ws.Cells.Item(37, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$cell.ClearFormats()
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '19.63'
$cell.ClearFormats()
$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +1.63%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(38, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Cronos'
$cell.ClearFormats()
$cell = $ws.Cells.Item(38, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell.ClearFormats()
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.105'
$cell.ClearFormats()
$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +22.15%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.113'
$cell.ClearFormats()
$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -4.81%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(40, 2)
$cell.NumberFormat = "@"
$cell.Value = 'USDe'
$cell.ClearFormats()
$cell = $ws.Cells.Item(40, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$cell.ClearFormats()
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.ClearFormats()
$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.08%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(41, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Aave'
$cell.ClearFormats()
$cell = $ws.Cells.Item(41, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell.ClearFormats()
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '180.33'
$cell.ClearFormats()
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.30%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(42, 2)
$cell.NumberFormat = "@"
$cell.Value = 'PolygonEcosystemToken'
$cell.ClearFormats()
$cell = $ws.Cells.Item(42, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$cell.ClearFormats()
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.345'
$cell.ClearFormats()
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +2.12%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(43, 2)
$cell.NumberFormat = "@"
$cell.Value = 'RenderToken'
$cell.ClearFormats()
$cell = $ws.Cells.Item(43, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$cell.ClearFormats()
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.98'
$cell.ClearFormats()
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.00%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(44, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Stacks'
$cell.ClearFormats()
$cell = $ws.Cells.Item(44, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell.ClearFormats()
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.65'
$cell.ClearFormats()
$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -1.43%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(45, 2)
$cell.NumberFormat = "@"
$cell.Value = 'OKB'
$cell.ClearFormats()
$cell = $ws.Cells.Item(45, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell.ClearFormats()
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '40.07'
$cell.ClearFormats()
$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.56%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(46, 2)
$cell.NumberFormat = "@"
$cell.Value = 'ImmutableX'
$cell.ClearFormats()
$cell = $ws.Cells.Item(46, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell.ClearFormats()
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.19'
$cell.ClearFormats()
$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -2.01%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(47, 2)
$cell.NumberFormat = "@"
$cell.Value = 'dogwifhat'
$cell.ClearFormats()
$cell = $ws.Cells.Item(47, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell.ClearFormats()
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.33'
$cell.ClearFormats()
$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.97%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(48, 2)
$cell.NumberFormat = "@"
$cell.Value = 'ARBITRUM'
$cell.ClearFormats()
$cell = $ws.Cells.Item(48, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell.ClearFormats()
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.571'
$cell.ClearFormats()
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +0.84%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(49, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Filecoin'
$cell.ClearFormats()
$cell = $ws.Cells.Item(49, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell.ClearFormats()
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.72'
$cell.ClearFormats()
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = '  -0.23%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(50, 2)
$cell.NumberFormat = "@"
$cell.Value = 'Mantle'
$cell.ClearFormats()
$cell = $ws.Cells.Item(50, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$cell.ClearFormats()
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.658'
$cell.ClearFormats()
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +7.16%  '
$cell.ClearFormats()
$cell = $ws.Cells.Item(51, 2)
$cell.NumberFormat = "@"
$cell.Value = 'InjectiveProtocol'
$cell.ClearFormats()
$cell = $ws.Cells.Item(51, 3)
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell.ClearFormats()
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '21.93'
$cell.ClearFormats()
$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = '  +5.49%  '
$cell.ClearFormats()
